$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated crypto price/volume data per the Sat Sep 30 2023 GitHub Actions refresh.
# Note: rows 17 and 18 swap content (WrappedBTC <-> Chainlink) plus their own value updates.
# Price column (D) values are plain text in this sheet; for values that look like pure
# numbers, force text formatting first (then restore the default style) so Excel keeps
# them as text instead of auto-converting to numeric cells.

$ws.Range("D2").Value = "27.024.89"
$ws.Range("E2").Value = "  +0.50%  "
$ws.Range("D3").Value = "1.675.86"
$ws.Range("E3").Value = "  +0.48%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "215.18"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.05%  "
$ws.Range("E6").Value = "  -1.02%  "
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("E8").Value = "  +2.31%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "21.40"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +5.40%  "
$ws.Range("E10").Value = "  -0.03%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0888"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.31%  "
$ws.Range("D12").Value = "1.913.22"
$ws.Range("E12").Value = "  +0.54%  "
$ws.Range("D13").Value = "1.677.73"
$ws.Range("E13").Value = "  +0.64%  "
$ws.Range("E15").Value = "  +1.42%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "66.31"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.50%  "
$ws.Range("B17").Value = "Chainlink"
$ws.Range("C17").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "8.21"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +3.10%  "
$ws.Range("B18").Value = "WrappedBTC"
$ws.Range("C18").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D18").Value = "27.032.88"
$ws.Range("E18").Value = "  +0.53%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "235.95"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.29%  "
$ws.Range("E20").Value = "  +0.67%  "
$ws.Range("E21").Value = "  -0.03%  "
$ws.Range("E22").Value = "  +1.34%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.23"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.98%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "148.10"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.02%  "
$ws.Range("E26").Value = "  +1.97%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "16.46"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +3.62%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.112"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.48%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.999"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.15%  "
$ws.Range("E30").Value = "  +0.45%  "
$ws.Range("E31").Value = "  -0.10%  "
$ws.Range("D33").Value = "1.542.32"
$ws.Range("E33").Value = "  +6.79%  "
$ws.Range("E34").Value = "  +1.21%  "
$ws.Range("E35").Value = "  +4.97%  "
$ws.Range("E36").Value = "  -1.19%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.590"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.80%  "
$ws.Range("E38").Value = "  +1.12%  "
$ws.Range("E39").Value = "  +1.97%  "
$ws.Range("E40").Value = "  +5.07%  "
$ws.Range("E41").Value = "  +0.00%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "67.82"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +2.81%  "
$ws.Range("E43").Value = "  -3.49%  "
$ws.Range("E44").Value = "  -2.09%  "
$ws.Range("D45").Value = "1.819.73"
$ws.Range("E45").Value = "  +0.67%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.779"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.61%  "
$ws.Range("E47").Value = "  -0.38%  "
$ws.Range("E48").Value = "  -0.08%  "
$ws.Range("E49").Value = "  +1.84%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.97"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +5.75%  "
$ws.Range("E51").Value = "  +0.31%  "
